$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.063.38'
$ws.Range('E2').Value = '  -0.59%  '
$ws.Range('D3').Value = '2.418.70'
$ws.Range('E3').Value = '  -1.14%  '
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').Value = '''552.91'
$ws.Range('E5').Value = '  -0.75%  '
$ws.Range('D6').Value = '''137.23'
$ws.Range('E6').Value = '  -1.43%  '
$ws.Range('D7').Value = '''0.999'
$ws.Range('E7').Value = '  -0.08%  '
$ws.Range('D8').Value = '''0.594'
$ws.Range('E8').Value = '  +3.71%  '
$ws.Range('E9').Value = '  -1.66%  '
$ws.Range('D10').Value = '''5.67'
$ws.Range('E10').Value = '  -2.31%  '
$ws.Range('E11').Value = '  -1.03%  '
$ws.Range('D12').Value = '''0.353'
$ws.Range('E12').Value = '  -2.53%  '
$ws.Range('D13').Value = '''25.19'
$ws.Range('E13').Value = '  +0.74%  '
$ws.Range('D14').Value = '2.847.61'
$ws.Range('E14').Value = '  -1.13%  '
$ws.Range('D15').Value = '59.944.47'
$ws.Range('E15').Value = '  -0.70%  '
$ws.Range('D16').Value = '''0.0000138'
$ws.Range('E16').Value = '  -2.12%  '
$ws.Range('D17').Value = '2.412.57'
$ws.Range('E17').Value = '  -1.20%  '
$ws.Range('D18').Value = '''11.29'
$ws.Range('E18').Value = '  -1.78%  '
$ws.Range('E19').Value = '  -0.67%  '
$ws.Range('D20').Value = '''328.02'
$ws.Range('E20').Value = '  -2.40%  '
$ws.Range('E21').Value = '  -3.45%  '
$ws.Range('D22').Value = '''1.00'
$ws.Range('E22').Value = '  +0.03%  '
$ws.Range('D23').Value = '''65.88'
$ws.Range('E23').Value = '  +1.80%  '
$ws.Range('D24').Value = '''0.176'
$ws.Range('E24').Value = '  +3.19%  '
$ws.Range('D25').Value = '''8.59'
$ws.Range('E25').Value = '  +0.53%  '
$ws.Range('E26').Value = '  +0.05%  '
$ws.Range('D27').Value = '''1.39'
$ws.Range('E27').Value = '  +0.63%  '
$ws.Range('D28').Value = '0.0₃0776'
$ws.Range('E28').Value = '  -2.75%  '
$ws.Range('E29').Value = '  -2.37%  '
$ws.Range('D30').Value = '''169.13'
$ws.Range('E30').Value = '  -1.04%  '
$ws.Range('E31').Value = '  -4.36%  '
$ws.Range('B32').Value = 'SuiNetwork'
$ws.Range('C32').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D32').Value = '''1.04'
$ws.Range('E32').Value = '  +0.78%  '
$ws.Range('B33').Value = 'EthereumClassic'
$ws.Range('C33').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D33').Value = '''18.57'
$ws.Range('E33').Value = '  -1.60%  '
$ws.Range('E34').Value = '  -0.02%  '
$ws.Range('E35').Value = '  -0.87%  '
$ws.Range('E36').Value = '  -0.04%  '
$ws.Range('D37').Value = '''4.19'
$ws.Range('E37').Value = '  -2.40%  '
$ws.Range('E38').Value = '  -2.08%  '
$ws.Range('D39').Value = '''323.82'
$ws.Range('E39').Value = '  +2.20%  '
$ws.Range('E40').Value = '  -3.52%  '
$ws.Range('D41').Value = '''3.66'
$ws.Range('E41').Value = '  -2.18%  '
$ws.Range('D42').Value = '''140.49'
$ws.Range('E42').Value = '  -2.64%  '
$ws.Range('E43').Value = '  +0.48%  '
$ws.Range('D44').Value = '''19.63'
$ws.Range('E44').Value = '  -1.62%  '
$ws.Range('D45').Value = '''0.0515'
$ws.Range('E45').Value = '  -2.00%  '
$ws.Range('D46').Value = '''0.577'
$ws.Range('E46').Value = '  +0.50%  '
$ws.Range('D47').Value = '''0.0224'
$ws.Range('E47').Value = '  -1.62%  '
$ws.Range('E48').Value = '  -5.46%  '
$ws.Range('E49').Value = '  +0.00%  '
$ws.Range('E50').Value = '  -5.03%  '
$ws.Range('E51').Value = '  -1.00%  '
